$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text for the latest version row (C8) to reflect
# the reverted "Plan" style replacement and restored 2D top-view hatching.
$ws.Range("C8").Value = "Annulation du remplacement de style en mode Plan (pour revenir à Dessin par défaut) et récupérer le hachurage en vue 2D de dessus"

# Move the sheet selection/scroll position: was C12, now C11 (view scrolled
# so row 7 is at the top).
$ws.Activate() | Out-Null
$ws.Range("A7").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C11").Select() | Out-Null
